# "added 4wk low sales check"
# Updates the rolling 4-week lookback derived columns (MyForecast, Inventory
# Coverage, Seasonality Index) on the "Forecast Comparison" sheet, and the
# dependent rollup totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: D (MyForecast), H (Inventory Coverage), L (Seasonality Index) ---

$wsForecast.Range("D2").Value  = 4
$wsForecast.Range("H2").Value  = 23.49
$wsForecast.Range("L2").Value  = 1.19

$wsForecast.Range("D3").Value  = 4
$wsForecast.Range("H3").Value  = 19.73
$wsForecast.Range("L3").Value  = 0.87

$wsForecast.Range("D4").Value  = 4
$wsForecast.Range("H4").Value  = 19.53
$wsForecast.Range("L4").Value  = 0.82

$wsForecast.Range("H5").Value  = 23.54
$wsForecast.Range("L5").Value  = 0.88

$wsForecast.Range("D6").Value  = 2
$wsForecast.Range("H6").Value  = 28.76
$wsForecast.Range("L6").Value  = 1.02

$wsForecast.Range("H7").Value  = 24.39
$wsForecast.Range("L7").Value  = 0.86

$wsForecast.Range("H8").Value  = 17.95
$wsForecast.Range("L8").Value  = 1.17

$wsForecast.Range("H9").Value  = 15.51
$wsForecast.Range("L9").Value  = 1.18

$wsForecast.Range("D10").Value = 4
$wsForecast.Range("H10").Value = 15.86
$wsForecast.Range("L10").Value = 1.13

$wsForecast.Range("H11").Value = 18.26
$wsForecast.Range("L11").Value = 0.95

$wsForecast.Range("H12").Value = 18.3
$wsForecast.Range("L12").Value = 0.99

$wsForecast.Range("D13").Value = 4
$wsForecast.Range("H13").Value = 13.93
$wsForecast.Range("L13").Value = 1.17

$wsForecast.Range("D14").Value = 4
$wsForecast.Range("H14").Value = 11.28
$wsForecast.Range("L14").Value = 0.98

$wsForecast.Range("D15").Value = 4
$wsForecast.Range("H15").Value = 10.28
$wsForecast.Range("L15").Value = 0.96

$wsForecast.Range("H16").Value = 11.18
$wsForecast.Range("L16").Value = 0.8100000000000001

$wsForecast.Range("H17").Value = 11.34
$wsForecast.Range("L17").Value = 1.08

# --- Summary: forecast rollups reflect the recomputed 4wk figures ---

$wsSummary.Range("B9").Value  = "65"
$wsSummary.Range("B10").Value = "33"
$wsSummary.Range("B11").Value = "18"
$wsSummary.Range("B12").Value = "5"
